$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 11
$ws.Range("D2").Value = 11

$ws.Range("C4").Value = 9
$ws.Range("D4").Value = 2

$ws.Range("D10").Select()
